$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.796.55"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.122.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.40"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.54"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.05%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.116.48"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +11.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.73"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.18"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.123"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.635.61"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.16"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.668.42"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.115.37"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "466.00"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.34"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.732"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.19"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.17"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.97"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +8.54%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.86"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0881"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +9.43%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.41%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.42"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +13.44%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.85"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "454.15"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.70"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.884.21"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.278"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.81"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.84"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.75"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.63%  "
